$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-12-20 07:08:42"

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("AA1")
    if ($cell.Text -eq "as_of_utc") {
        $ws.Range("AA2:AA26").Value = $newTimestamp
    }
}
